$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-03-03 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-04 Tuesday", 2) | Out-Null

# Update each division-problem cell in the table, addressed by (row, column)
# so that the duplicate "82÷6=13, 4" source values are each replaced with the
# correct, position-specific target value.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "93÷9=10, 3"
$t.Cell(1, 2).Range.Text = "59÷4=14, 3"
$t.Cell(1, 3).Range.Text = "21÷3=7, 0"
$t.Cell(1, 4).Range.Text = "23÷4=5, 3"
$t.Cell(1, 5).Range.Text = "87÷2=43, 1"

$t.Cell(5, 1).Range.Text = "64÷9=7, 1"
$t.Cell(5, 2).Range.Text = "71÷8=8, 7"
$t.Cell(5, 3).Range.Text = "63÷5=12, 3"
$t.Cell(5, 4).Range.Text = "30÷4=7, 2"
$t.Cell(5, 5).Range.Text = "73÷7=10, 3"

$t.Cell(9, 1).Range.Text = "57÷9=6, 3"
$t.Cell(9, 2).Range.Text = "24÷9=2, 6"
$t.Cell(9, 3).Range.Text = "27÷9=3, 0"
$t.Cell(9, 4).Range.Text = "86÷6=14, 2"
$t.Cell(9, 5).Range.Text = "23÷9=2, 5"

$t.Cell(13, 1).Range.Text = "34÷3=11, 1"
$t.Cell(13, 2).Range.Text = "46÷6=7, 4"
$t.Cell(13, 3).Range.Text = "25÷2=12, 1"
$t.Cell(13, 4).Range.Text = "41÷4=10, 1"
$t.Cell(13, 5).Range.Text = "37÷9=4, 1"

$t.Cell(17, 1).Range.Text = "20÷2=10, 0"
$t.Cell(17, 2).Range.Text = "98÷4=24, 2"
$t.Cell(17, 3).Range.Text = "43÷7=6, 1"
$t.Cell(17, 4).Range.Text = "82÷6=13, 4"
$t.Cell(17, 5).Range.Text = "29÷5=5, 4"

